# Applies the updated cryptos list values (prices / 1h volume %) and
# the row-50/row-51 coin swap (HuobiToken <-> NEARProtocol), matching the
# commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "44.030.66"
Set-TextCell "E2" "  +4.32%  "

Set-TextCell "D3" "2.256.78"
Set-TextCell "E3" "  +1.95%  "

Set-TextCell "E4" "  -0.05%  "

Set-TextCell "D5" "229.47"
Set-TextCell "E5" "  -0.52%  "

Set-TextCell "D6" "0.632"
Set-TextCell "E6" "  +2.36%  "

Set-TextCell "D7" "63.32"
Set-TextCell "E7" "  +4.36%  "

Set-TextCell "E8" "  -0.10%  "

Set-TextCell "D9" "0.442"
Set-TextCell "E9" "  +10.23%  "

Set-TextCell "D10" "0.102"
Set-TextCell "E10" "  +13.35%  "

Set-TextCell "E11" "  -0.38%  "

Set-TextCell "D12" "25.92"
Set-TextCell "E12" "  +17.39%  "

Set-TextCell "D13" "0.106"
Set-TextCell "E13" "  +2.14%  "

Set-TextCell "D14" "2.591.60"
Set-TextCell "E14" "  +1.70%  "

Set-TextCell "D15" "15.56"
Set-TextCell "E15" "  +0.65%  "

Set-TextCell "D16" "6.12"
Set-TextCell "E16" "  +9.84%  "

Set-TextCell "D17" "0.839"
Set-TextCell "E17" "  +5.61%  "

Set-TextCell "D18" "2.270.62"
Set-TextCell "E18" "  +2.00%  "

Set-TextCell "D19" "43.935.48"
Set-TextCell "E19" "  +4.24%  "

Set-TextCell "E20" "  +7.75%  "

Set-TextCell "D21" "73.07"
Set-TextCell "E21" "  +1.46%  "

Set-TextCell "E22" "  -2.97%  "

Set-TextCell "D23" "251.28"
Set-TextCell "E23" "  +3.21%  "

Set-TextCell "E24" "  +0.01%  "

Set-TextCell "E25" "  +0.68%  "

Set-TextCell "E26" "  -2.77%  "

Set-TextCell "D27" "9.98"
Set-TextCell "E27" "  +4.02%  "

Set-TextCell "D28" "3.25"
Set-TextCell "E28" "  +22.66%  "

Set-TextCell "D29" "171.83"
Set-TextCell "E29" "  +1.41%  "

Set-TextCell "E30" "  +2.26%  "

Set-TextCell "D31" "0.136"
Set-TextCell "E31" "  -3.21%  "

Set-TextCell "D32" "1.38"
Set-TextCell "E32" "  -4.69%  "

Set-TextCell "E33" "  +2.73%  "

Set-TextCell "D34" "0.0684"
Set-TextCell "E34" "  +5.33%  "

Set-TextCell "D35" "4.71"
Set-TextCell "E35" "  +2.10%  "

Set-TextCell "D36" "4.84"
Set-TextCell "E36" "  -2.47%  "

Set-TextCell "D37" "3.81"
Set-TextCell "E37" "  +7.20%  "

Set-TextCell "D38" "6.57"
Set-TextCell "E38" "  +3.75%  "

Set-TextCell "E39" "  -1.73%  "

Set-TextCell "E40" "  +2.80%  "

Set-TextCell "E41" "  -0.05%  "

Set-TextCell "D42" "17.32"
Set-TextCell "E42" "  +8.40%  "

Set-TextCell "D43" "8.19"
Set-TextCell "E43" "  -3.98%  "

Set-TextCell "D44" "0.0962"
Set-TextCell "E44" "  +0.38%  "

Set-TextCell "D45" "97.17"
Set-TextCell "E45" "  +0.41%  "

Set-TextCell "E46" "  -0.63%  "

Set-TextCell "D47" "0.000210"
Set-TextCell "E47" "  -9.28%  "

Set-TextCell "D48" "4.34"
Set-TextCell "E48" "  -0.53%  "

Set-TextCell "D49" "1.432.88"
Set-TextCell "E49" "  -1.47%  "

Set-TextCell "B50" "NEARProtocol"
Set-TextCell "C50" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D50" "2.27"
Set-TextCell "E50" "  +2.89%  "

Set-TextCell "B51" "HuobiToken"
Set-TextCell "C51" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D51" "2.75"
Set-TextCell "E51" "  +0.48%  "
